$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "43.398.68"
$ws.Range("E2").Value = "  -0.96%  "

# Row 3
$ws.Range("D3").Value = "2.333.14"
$ws.Range("E3").Value = "  -0.73%  "

# Row 4
$ws.Range("E4").Value = "  -0.12%  "

# Row 5
$ws.Range("D5").Value = "'239.08"
$ws.Range("E5").Value = "  -0.26%  "

# Row 6
$ws.Range("D6").Value = "'0.659"
$ws.Range("E6").Value = "  -2.37%  "

# Row 7
$ws.Range("D7").Value = "'73.38"
$ws.Range("E7").Value = "  -0.60%  "

# Row 8
$ws.Range("E8").Value = "  +0.06%  "

# Row 9
$ws.Range("D9").Value = "'0.590"
$ws.Range("E9").Value = "  -2.38%  "

# Row 10
$ws.Range("E10").Value = "  +1.97%  "

# Row 11
$ws.Range("D11").Value = "'59.37"
$ws.Range("E11").Value = "  +2.37%  "

# Row 12
$ws.Range("D12").Value = "'36.14"
$ws.Range("E12").Value = "  +9.74%  "

# Row 13
$ws.Range("E13").Value = "  +0.66%  "

# Row 14
$ws.Range("D14").Value = "'7.24"
$ws.Range("E14").Value = "  -1.25%  "

# Row 15
$ws.Range("D15").Value = "2.680.72"
$ws.Range("E15").Value = "  -0.88%  "

# Row 16
$ws.Range("D16").Value = "'16.20"
$ws.Range("E16").Value = "  -1.55%  "

# Row 17
$ws.Range("D17").Value = "'0.923"
$ws.Range("E17").Value = "  +2.12%  "

# Row 18
$ws.Range("D18").Value = "2.325.51"
$ws.Range("E18").Value = "  -1.27%  "

# Row 19
$ws.Range("D19").Value = "43.345.09"
$ws.Range("E19").Value = "  -0.94%  "

# Row 20
$ws.Range("D20").Value = "'0.0000103"
$ws.Range("E20").Value = "  +1.45%  "

# Row 21
$ws.Range("D21").Value = "'6.52"
$ws.Range("E21").Value = "  -3.83%  "

# Row 22
$ws.Range("D22").Value = "'76.68"
$ws.Range("E22").Value = "  -0.91%  "

# Row 23
$ws.Range("D23").Value = "'250.76"
$ws.Range("E23").Value = "  -2.53%  "

# Row 24
$ws.Range("E24").Value = "  +0.05%  "

# Row 25
$ws.Range("D25").Value = "'3.75"
$ws.Range("E25").Value = "  +0.46%  "

# Row 26
$ws.Range("E26").Value = "  -8.22%  "

# Row 27
$ws.Range("D27").Value = "'2.48"
$ws.Range("E27").Value = "  -0.20%  "

# Row 28
$ws.Range("D28").Value = "'10.55"
$ws.Range("E28").Value = "  -0.88%  "

# Row 29
$ws.Range("D29").Value = "'2.28"
$ws.Range("E29").Value = "  +0.93%  "

# Row 30
$ws.Range("D30").Value = "'173.65"
$ws.Range("E30").Value = "  -1.53%  "

# Row 31
$ws.Range("D31").Value = "'21.90"
$ws.Range("E31").Value = "  -3.94%  "

# Row 32
$ws.Range("D32").Value = "'0.127"
$ws.Range("E32").Value = "  -2.31%  "

# Row 33
$ws.Range("D33").Value = "'0.133"
$ws.Range("E33").Value = "  -2.95%  "

# Row 34
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").Value = "'5.56"
$ws.Range("E34").Value = "  +1.38%  "

# Row 35
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "'0.0744"
$ws.Range("E35").Value = "  -2.53%  "

# Row 36
$ws.Range("D36").Value = "'5.10"
$ws.Range("E36").Value = "  -1.62%  "

# Row 37
$ws.Range("D37").Value = "'3.75"
$ws.Range("E37").Value = "  -0.63%  "

# Row 38
$ws.Range("D38").Value = "'6.49"
$ws.Range("E38").Value = "  +3.61%  "

# Row 39
$ws.Range("D39").Value = "'2.36"
$ws.Range("E39").Value = "  +0.19%  "

# Row 40
$ws.Range("D40").Value = "'0.0278"
$ws.Range("E40").Value = "  -0.89%  "

# Row 41
$ws.Range("D41").Value = "'20.99"
$ws.Range("E41").Value = "  +10.76%  "

# Row 42
$ws.Range("B42").Value = "MultiversX"
$ws.Range("C42").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D42").Value = "'67.44"
$ws.Range("E42").Value = "  +3.36%  "

# Row 43
$ws.Range("B43").Value = "FTXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D43").Value = "'5.28"
$ws.Range("E43").Value = "  +12.00%  "

# Row 44
$ws.Range("D44").Value = "'0.107"
$ws.Range("E44").Value = "  -4.79%  "

# Row 45
$ws.Range("D45").Value = "'9.11"
$ws.Range("E45").Value = "  +0.39%  "

# Row 46
$ws.Range("D46").Value = "'0.198"
$ws.Range("E46").Value = "  -3.67%  "

# Row 47
$ws.Range("D47").Value = "'2.50"
$ws.Range("E47").Value = "  -0.05%  "

# Row 48
$ws.Range("B48").Value = "BinanceUSD"
$ws.Range("C48").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D48").Value = "'1.00"
$ws.Range("E48").Value = "  -0.19%  "

# Row 49
$ws.Range("B49").Value = "TrustWalletToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D49").Value = "'1.23"
$ws.Range("E49").Value = "  -0.28%  "

# Row 50
$ws.Range("B50").Value = "SynthetixNetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D50").Value = "'4.48"
$ws.Range("E50").Value = "  +12.43%  "

# Row 51
$ws.Range("B51").Value = "ARBITRUM"
$ws.Range("C51").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D51").Value = "'1.15"
$ws.Range("E51").Value = "  -1.08%  "
